# Apply the updated symbol-list snapshot values (price, 1h volume %, hour)
# to Sheet1. Target cells store these as literal text (not numbers), matching
# the original inline-string cell type, so each value is written with a
# leading apostrophe (forces Excel to keep it as text) and the cell style is
# reset to "Normal" afterward so no stray number-format/quote-prefix style
# sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.33%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'16"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'31.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.25%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'16"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.080"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.29%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'16"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08035"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'9.01%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'16"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'2.705"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'47.83%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'16"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'7.807"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.89%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'16"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'3.792"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.57%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'16"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.9271"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.08%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'16"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1746"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.15%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'16"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.07197"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.01%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'16"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.08969"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'10.52%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'16"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.03017"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.94%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'16"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.09994"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.53%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'16"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.001489"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.44%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'16"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.005939"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.38%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'16"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'3.535"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'16"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'2.246"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.03%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'16"
$ws.Range("G18").Style = "Normal"
$ws.Range("E19").Value = "'-0.89%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'16"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.86%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'16"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'4.188"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-9.28%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'16"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.1645"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.21%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'16"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04584"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.56%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'16"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.001239"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.74%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'16"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004424"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.48%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'16"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001196"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.79%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'16"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003421"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'82.76%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'16"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'16"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'16"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'16"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'16"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'16"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'16"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'16"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'16"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'16"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'16"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'16"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.01768"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.01%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'16"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04485"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'16"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006842"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'16"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1348"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.13%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'16"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.002128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.74%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'16"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.009833"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-9.26%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'16"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.72%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'16"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000747"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.39%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'16"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.008733"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-14.50%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'16"
$ws.Range("G47").Style = "Normal"
$ws.Range("E48").Value = "'11.23%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'16"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002092"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.39%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'16"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.32%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'16"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'16"
$ws.Range("G51").Style = "Normal"
